$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.91
$ws.Range("C3").Value = 2.0988
$ws.Range("C4").Value = 0.7076
$ws.Range("C5").Value = 0.1558
$ws.Range("C6").Value = 0.1409
$ws.Range("C7").Value = 0.3932
